$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Counts (cps)" (column B) and "Error (cps)" (column C) values
# for rows 2-19, per the draft results rewrite.

$values = @{
    2  = @(1.27366666666667, 0.09216252)
    3  = @(1.74133333333333, 0.107753706666667)
    4  = @(1.251, 0.09139806)
    5  = @(5.52033333333333, 0.191886786666667)
    6  = @(4.83366666666667, 0.17952238)
    7  = @(5.088, 0.1841856)
    8  = @(18.8956666666667, 0.35486062)
    9  = @(19.5463333333333, 0.360825313333333)
    10 = @(18.6393333333333, 0.352656186666667)
    11 = @(37.508, 0.49960656)
    12 = @(38.086, 0.50349692)
    13 = @(38.2003333333333, 0.505008406666667)
    14 = @(90.1833333333333, 0.777380333333333)
    15 = @(93.4766666666667, 0.788943066666667)
    16 = @(98.6423333333333, 0.81083998)
    17 = @(191.358, 1.1290122)
    18 = @(189.927666666667, 1.12817034)
    19 = @(190.764666666667, 1.12932682666667)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
}
